$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10 (Leve Item ID 1959)
$ws.Range("H10").Value = 876
$ws.Range("I10").Value = 749.5
$ws.Range("J10").Value = 1002.5
$ws.Range("K10").Value = 749.5
$ws.Range("L10").Value = 1002.5
$ws.Range("M10").Value = -456.5
$ws.Range("N10").Value = -1588.5

# Row 19 (Leve Item ID 7015)
$ws.Range("H19").Value = 1035.5
$ws.Range("I19").Value = 1215
$ws.Range("K19").Value = 1215
$ws.Range("M19").Value = -1040

# Row 100 (Leve Item ID 19906)
$ws.Range("H100").Value = 3167.5
$ws.Range("I100").Value = 1516.3636
$ws.Range("K100").Value = 1516.3636
$ws.Range("M100").Value = -975.3635999999999

# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 2000
$ws.Range("I106").Value = 2000
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2000
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1369
$ws.Range("N106").Value = $null

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 4637.75
$ws.Range("I116").Value = 4100.8335
$ws.Range("K116").Value = 4100.8335
$ws.Range("M116").Value = -658.8334999999997

$ws = $wb.Worksheets.Item("ARM")
# Row 97 (Leve Item ID 19941)
$ws.Range("H97").Value = 917.2727
$ws.Range("I97").Value = 776.5714
$ws.Range("K97").Value = 776.5714
$ws.Range("M97").Value = -280.5714

$ws = $wb.Worksheets.Item("BSM")
# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 3600
$ws.Range("I99").Value = 950
$ws.Range("K99").Value = 950
$ws.Range("M99").Value = 548

# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 2594.818
$ws.Range("I105").Value = 2694.3
$ws.Range("K105").Value = 2694.3
$ws.Range("M105").Value = -947.3000000000002

# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 8349.286
$ws.Range("I107").Value = 7784.4614
$ws.Range("K107").Value = 7784.4614
$ws.Range("M107").Value = -5864.4614

# Row 109 (Leve Item ID 27096)
$ws.Range("H109").Value = 40000
$ws.Range("J109").Value = 40000
$ws.Range("L109").Value = 40000
$ws.Range("N109").Value = -42774

$ws = $wb.Worksheets.Item("CRP")
# Row 22 (Leve Item ID 5367)
$ws.Range("H22").Value = 2200.923
$ws.Range("I22").Value = 1200.25
$ws.Range("J22").Value = 2645.6667
$ws.Range("K22").Value = 1200.25
$ws.Range("L22").Value = 2645.6667
$ws.Range("M22").Value = -850.25
$ws.Range("N22").Value = -3345.6667

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 1605.9333
$ws.Range("I58").Value = 1601.5454
$ws.Range("K58").Value = 1601.5454
$ws.Range("M58").Value = -1398.5454

# Row 70 (Leve Item ID 12011)
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null

# Row 73 (Leve Item ID 12011)
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null

# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 7600
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").Value = $null

# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 7600
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").Value = $null

# Row 94 (Leve Item ID 32934)
$ws.Range("H94").Value = 5499
$ws.Range("J94").Value = 5499
$ws.Range("L94").Value = 5499
$ws.Range("N94").Value = -6401

# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 765.05554
$ws.Range("I107").Value = 497.1
$ws.Range("K107").Value = 497.1
$ws.Range("M107").Value = 1422.9

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 1605.9333
$ws.Range("I136").Value = 1601.5454
$ws.Range("K136").Value = 4804.6362
$ws.Range("M136").Value = -2254.6362

$ws = $wb.Worksheets.Item("CUL")
# Row 7 (Leve Item ID 4728)
$ws.Range("H7").Value = 655.4167
$ws.Range("I7").Value = 642.2727
$ws.Range("K7").Value = 1926.8181
$ws.Range("M7").Value = -1814.8181

# Row 18 (Leve Item ID 36056)
$ws.Range("H18").Value = 4889.7
$ws.Range("I18").Value = 3574.5
$ws.Range("K18").Value = 10723.5
$ws.Range("M18").Value = -10554.5

# Row 19 (Leve Item ID 4682)
$ws.Range("H19").Value = 2120.4443
$ws.Range("I19").Value = 2120.4443
$ws.Range("K19").Value = 6361.3329
$ws.Range("M19").Value = -6187.3329

# Row 22 (Leve Item ID 4697)
$ws.Range("H22").Value = 763.5
$ws.Range("J22").Value = 763.5
$ws.Range("L22").Value = 2290.5
$ws.Range("N22").Value = -2628.5

# Row 25 (Leve Item ID 4709)
$ws.Range("H25").Value = 2874.75
$ws.Range("I25").Value = 499.66666
$ws.Range("J25").Value = 10000
$ws.Range("K25").Value = 1498.99998
$ws.Range("L25").Value = 30000
$ws.Range("M25").Value = -1329.99998
$ws.Range("N25").Value = -30338

# Row 27 (Leve Item ID 4697)
$ws.Range("H27").Value = 763.5
$ws.Range("J27").Value = 763.5
$ws.Range("L27").Value = 2290.5
$ws.Range("N27").Value = -2494.5

# Row 29 (Leve Item ID 4698)
$ws.Range("H29").Value = 292.25
$ws.Range("I29").Value = 222
$ws.Range("K29").Value = 666
$ws.Range("M29").Value = -389

# Row 30 (Leve Item ID 4709)
$ws.Range("H30").Value = 2874.75
$ws.Range("I30").Value = 499.66666
$ws.Range("J30").Value = 10000
$ws.Range("K30").Value = 1498.99998
$ws.Range("L30").Value = 30000
$ws.Range("M30").Value = -1396.99998
$ws.Range("N30").Value = -30204

# Row 36 (Leve Item ID 4732)
$ws.Range("H36").Value = 646
$ws.Range("I36").Value = 475.2
$ws.Range("J36").Value = 1500
$ws.Range("K36").Value = 1425.6
$ws.Range("L36").Value = 4500
$ws.Range("M36").Value = -1256.6
$ws.Range("N36").Value = -4838

# Row 94 (Leve Item ID 19811)
$ws.Range("H94").Value = 7333.25
$ws.Range("I94").Value = 3200
$ws.Range("J94").Value = 10285.571
$ws.Range("K94").Value = 9600
$ws.Range("L94").Value = 30856.713
$ws.Range("M94").Value = -8924
$ws.Range("N94").Value = -32208.713

# Row 106 (Leve Item ID 19819)
$ws.Range("H106").Value = 10433.333
$ws.Range("J106").Value = 15000
$ws.Range("L106").Value = 45000
$ws.Range("N106").Value = -46892

# Row 132 (Leve Item ID 43972)
$ws.Range("H132").Value = 1269.5714
$ws.Range("I132").Value = 998
$ws.Range("K132").Value = 8982
$ws.Range("M132").Value = -6452

# Row 141 (Leve Item ID 44076)
$ws.Range("H141").Value = 4988.6
$ws.Range("I141").Value = 2121.5
$ws.Range("J141").Value = 6900
$ws.Range("K141").Value = 6364.5
$ws.Range("L141").Value = 20700
$ws.Range("M141").Value = -1184.5
$ws.Range("N141").Value = -31060

$ws = $wb.Worksheets.Item("LTW")
# Row 40 (Leve Item ID 36248)
$ws.Range("H40").Value = 5433
$ws.Range("I40").Value = 4154.727
$ws.Range("K40").Value = 4154.727
$ws.Range("M40").Value = -4018.727

# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 3099.5625
$ws.Range("I46").Value = 897
$ws.Range("J46").Value = 3414.2144
$ws.Range("K46").Value = 897
$ws.Range("L46").Value = 3414.2144
$ws.Range("M46").Value = -709
$ws.Range("N46").Value = -3790.2144

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 5577.1665
$ws.Range("I93").Value = 3401.375
$ws.Range("J93").Value = 7317.8
$ws.Range("K93").Value = 3401.375
$ws.Range("L93").Value = 7317.8
$ws.Range("M93").Value = -2153.375
$ws.Range("N93").Value = -9813.799999999999

$ws = $wb.Worksheets.Item("WVR")
# Row 70 (Leve Item ID 11979)
$ws.Range("H70").Value = 142896420
$ws.Range("J70").Value = 166705000
$ws.Range("L70").Value = 166705000
$ws.Range("N70").Value = -166705630

# Row 73 (Leve Item ID 11979)
$ws.Range("H73").Value = 142896420
$ws.Range("J73").Value = 166705000
$ws.Range("L73").Value = 166705000
$ws.Range("N73").Value = -166707184

# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 1627.4286
$ws.Range("J107").Value = 1997.5
$ws.Range("L107").Value = 1997.5
$ws.Range("N107").Value = -9832.5
